# Badge Order Tracking.xlsx - "Add files via upload" edit
# Adds a "Duty" line (LCSC component duty), a "Rotary Knobs" line, and a
# "Carl Total" formula row to the bottom of the tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20: record the arrival date for the 14500 Batteries line ---
$ws.Range("H20").Value = 44820
$ws.Range("H20").NumberFormat = "d-mmm"

# --- Row 26 (new): Duty on LCSC Component Order ---
$ws.Range("A26").Value = "Duty"
$ws.Range("B26").Value = "Duty on LCSC Component Order"
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = "DHL"
$ws.Range("E26").Value = 380
$ws.Range("F26").Value = "Carl"
$ws.Range("H26").Value = 44819
$ws.Range("H26").NumberFormat = "d-mmm"
$ws.Range("I26").Value = 44819
$ws.Range("I26").NumberFormat = "d-mmm"

# --- Row 27 (previously blank): Rotary Knobs ---
$ws.Range("A27").Value = "Components"
$ws.Range("B27").Value = "Rotary Knobs"
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = "AliExpress"
$ws.Range("E27").Value = 323.49
$ws.Range("F27").Value = "Carl"
$ws.Range("H27").Value = 44824
$ws.Range("H27").NumberFormat = "d-mmm"

# --- Row 38 (new): Carl Total ---
$ws.Range("E38").Value = "Carl Total"
$ws.Range("F38").Formula = "=E18+E19+E20+E21+E26+E27"
